# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Latest Handoff Datetime" on the per-locale sheets for the last tracked
# file (ead9b05c-6df6-4f2d-9561-2cf7d1a36e36.md), reflecting a fresh handoff
# xliff-generation pass.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-10-18 02:48:23"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-10-18 02:48:00"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-10-18 02:48:23"
